# edit.ps1
# Applies:
#  1. Update the cached "datetimeFigureOut" date field text
#     (2/11/25 -> 2/17/25) on the slide master and every slide layout's
#     "Date Placeholder" shape.
#  2. On slide 4 ("To submit your homework..."):
#       - "homework/hw5 folder" -> "homework/hw4 folder"
#       - collapse the "- change category to [ HW 5 ] " run-triplet into
#         a single run reading "- change category to [ HW 4 ] "

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text: 2/11/25 -> 2/17/25 (slide master + layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $full = $tf.TextRange
                if ($full.Text -eq "2/11/25") {
                    $full.Text = "2/17/25"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2. Slide 4 content edits
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$content = $slide4.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

# 2a. "homework/hw5 folder" -> "homework/hw4 folder" (run 1 of paragraph 2)
$para2 = $tr.Paragraphs(2, 1)
$run1Len = "1. Save your data visualization to the homework/hw5 folder using [".Length
$run1 = $para2.Characters(1, $run1Len)
if ($run1.Text -eq "1. Save your data visualization to the homework/hw5 folder using [") {
    $run1.Text = "1. Save your data visualization to the homework/hw4 folder using ["
}

# 2b. Collapse "	- change category to [ " / "HW 5 " / "] " into one run
$para5 = $tr.Paragraphs(5, 1)
# Paragraphs()/.Text includes the trailing paragraph-mark (CR) character,
# so compare against the content length (Length - 1), not the raw .Text.
$para5Content = $para5.Characters(1, $para5.Length - 1)
if ($para5Content.Text -eq "`t- change category to [ HW 5 ] ") {
    $lead = "`t- change category to [ ".Length
    $tail = $para5.Characters($lead + 1, $para5.Length - 1 - $lead)
    $tail.Delete()
    $kept = $tr.Paragraphs(5, 1)
    $keptContent = $kept.Characters(1, $kept.Length - 1)
    $keptContent.Text = "`t- change category to [ HW 4 ] "
}
